# Applies the cryptos.xlsx "symbol list" data refresh described by the commit:
# "Updated symbol list on Thu Jan  5 05:32:18 UTC 2023 with GitHub Actions"
#
# Every target cell holds its value as literal text (inlineStr/shared string),
# not a number or percentage, matching the upstream data feed format. Assigning
# a numeric- or percent-looking string straight to Range.Value makes Excel auto-
# convert it to a real number (and stamp on a % NumberFormat), so every write is
# led with an apostrophe (Excel's "treat as text" quote-prefix) and the style is
# then snapped back to Normal so no stray formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue "D2" '258.71'
Set-TextValue "E2" '1.23%'
Set-TextValue "D3" '26.73'
Set-TextValue "E3" '-4.94%'
Set-TextValue "D4" '4.830'
Set-TextValue "E4" '-9.79%'
Set-TextValue "D5" '0.05971'
Set-TextValue "E5" '2.48%'
Set-TextValue "D6" '6.687'
Set-TextValue "E6" '-0.44%'
Set-TextValue "E7" '1.20%'
Set-TextValue "D8" '0.9532'
Set-TextValue "E8" '5.02%'
Set-TextValue "D9" '0.1417'
Set-TextValue "E9" '-0.37%'
Set-TextValue "D10" '0.03589'
Set-TextValue "E10" '3.77%'
Set-TextValue "D11" '0.07224'
Set-TextValue "E11" '0.53%'
Set-TextValue "D12" '0.03144'
Set-TextValue "E12" '-1.27%'
Set-TextValue "D13" '0.09238'
Set-TextValue "E13" '-0.07%'
Set-TextValue "D14" '0.001539'
Set-TextValue "E14" '-0.06%'
Set-TextValue "B15" 'One'
Set-TextValue "C15" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue "D15" '0.0006077'
Set-TextValue "E15" '0.46%'
Set-TextValue "B16" 'TigerCash'
Set-TextValue "C16" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D16" '0.006084'
Set-TextValue "E16" '2.52%'
Set-TextValue "B17" 'LEO'
Set-TextValue "C17" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.486'
Set-TextValue "E17" '-0.31%'
Set-TextValue "B18" 'GateToken'
Set-TextValue "C18" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D18" '3.223'
Set-TextValue "E18" '-0.05%'
Set-TextValue "B19" 'BTSEToken'
Set-TextValue "C19" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D19" '2.239'
Set-TextValue "E19" '-1.46%'
Set-TextValue "D20" '0.3136'
Set-TextValue "E20" '-0.97%'
Set-TextValue "E21" '-2.12%'
Set-TextValue "E22" '-0.25%'
Set-TextValue "D23" '0.04229'
Set-TextValue "E23" '1.71%'
Set-TextValue "E24" '0.03%'
Set-TextValue "E25" '-0.37%'
Set-TextValue "D26" '0.004515'
Set-TextValue "E26" '-7.11%'
Set-TextValue "E27" '0.01%'
Set-TextValue "E28" '-23.00%'
Set-TextValue "D40" '0.03848'
Set-TextValue "E40" '0.05%'
Set-TextValue "D41" '0.005958'
Set-TextValue "E41" '3.69%'
Set-TextValue "D42" '0.1104'
Set-TextValue "E42" '0.48%'
Set-TextValue "D43" '0.002299'
Set-TextValue "E43" '4.55%'
Set-TextValue "D44" '0.01048'
Set-TextValue "E44" '6.28%'
Set-TextValue "E45" '3.73%'
Set-TextValue "E46" '-0.01%'
Set-TextValue "D47" '0.1090'
Set-TextValue "E47" '8.93%'
Set-TextValue "E48" '-3.67%'
Set-TextValue "E49" '-0.01%'
Set-TextValue "E50" '-0.01%'
